$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-02-08 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-09 Sunday", 2)

# Update the division problems in the table, addressed by (row, column)
# so that duplicate values introduced mid-edit don't get clobbered by
# later replacements.
$tbl = $d.Tables.Item(1)

$cellValues = @{
    1  = @("80÷7=", "45÷4=", "62÷8=", "47÷6=", "89÷6=")
    5  = @("24÷5=", "65÷3=", "29÷6=", "59÷6=", "64÷7=")
    9  = @("29÷6=", "30÷4=", "16÷7=", "33÷7=", "33÷6=")
    13 = @("14÷5=", "25÷3=", "39÷8=", "72÷3=", "10÷6=")
    17 = @("40÷4=", "15÷3=", "82÷3=", "94÷6=", "36÷7=")
}

foreach ($rowIndex in $cellValues.Keys) {
    $values = $cellValues[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cellRange = $cell.Range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $values[$col - 1]
    }
}
